$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 26 corresponds to "Open year" 2024.
# Update Energy Storage (column C) value from 1 to 3
$ws.Range("C26").Value = 3

# Update Solar (column E) value from 4 to 8
$ws.Range("E26").Value = 8
